$d = $word.ActiveDocument

# Locate the "Discussion" paragraph that starts with "First, in every fourth
# quarter ..." (the one right after the e-payment paragraph and right before
# the "E-payment" bulleted list) and rewrite its contents in one shot so the
# resulting run layout matches the authored edit exactly:
#   - "word" -> "phrase" in "explain the meaning of the word/phrase
#     "holiday consumption"" (kept in the existing red run)
#   - "festival" -> "festivals"
#   - "triggers" -> "trigger"
#   - "to buy products to treat themselves." -> "to buy products as treating
#     themselves. " (extra trailing space)
#   - the paragraph mark itself picks up red run-formatting (w:pPr/w:rPr),
#     matching what Word leaves behind after this kind of in-place edit.
$targetParagraph = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "First, in every fourth quarter*") {
        $targetParagraph = $p
        break
    }
}

if ($targetParagraph -eq $null) {
    throw "Could not locate the target paragraph"
}

$newParagraphXml = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:rPr><w:color w:val="FF0000"/></w:rPr></w:pPr><w:r><w:t>First, in every fourth quarter of the year, both Alibaba and Amazon’s online sales revenue increased</w:t></w:r><w:r><w:t xml:space="preserve"> to the top of the year. This phenomenon is usually come with holiday consumption, </w:t></w:r><w:r><w:rPr><w:color w:val="FF0000"/></w:rPr><w:t xml:space="preserve">explain the meaning of the </w:t></w:r><w:r><w:rPr><w:color w:val="FF0000"/></w:rPr><w:t>phrase</w:t></w:r><w:r><w:rPr><w:color w:val="FF0000"/></w:rPr><w:t xml:space="preserve"> “holiday consumption”</w:t></w:r><w:r><w:t>. Significant festival</w:t></w:r><w:r><w:t>s</w:t></w:r><w:r><w:t xml:space="preserve"> such as Christmas</w:t></w:r><w:r><w:t xml:space="preserve"> and the eleventh of November are highly trigger c</w:t></w:r><w:r><w:t>onsumer</w:t></w:r><w:r><w:t>’s</w:t></w:r><w:r><w:t xml:space="preserve"> desire</w:t></w:r><w:r><w:t xml:space="preserve"> to buy products </w:t></w:r><w:r><w:t xml:space="preserve">as </w:t></w:r><w:r><w:t>treat</w:t></w:r><w:r><w:t>ing</w:t></w:r><w:r><w:t xml:space="preserve"> themselves.</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@

$targetParagraph.Range.InsertXML($newParagraphXml)
